$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-72 down to 23-73.
$ws.Rows(22).Insert()

# Populate the newly inserted row 22 with the new weekly data entry.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44987
$ws.Range("D22").NumberFormat = $ws.Range("D23").NumberFormat
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100107
$ws.Range("H22").Value = "Otros"
$ws.Range("I22").Value = 100107011
$ws.Range("J22").Value = "Tuna"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 40
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("Q22").Value = "$/caja 16 kilos"
$ws.Range("R22").Value = "Provincia de Los Andes"
$ws.Range("S22").Value = 938
$ws.Range("T22").Value = 16
